{"js": "// Remove the \"Step 2:\" paragraph and merge the trailing (empty, bookmark-only)\n// paragraph into the preceding \"composer require intervention/image\" paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the paragraphs we need by their text content.\nlet stepTwoPara = null;\nlet composerPara = null;\nconst items = paragraphs.items;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (t.trim() === \"Step 2:\") {\n    stepTwoPara = items[i];\n  } else if (t.trim() === \"composer require intervention/image\") {\n    composerPara = items[i];\n  }\n}\n\nif (!stepTwoPara || !composerPara) {\n  throw new Error(\"Could not locate the expected paragraphs (Step 2: / composer require...).\");\n}\n\n// Delete the \"Step 2:\" paragraph entirely (text + its paragraph mark).\nstepTwoPara.delete();\nawait context.sync();\n\n// After removing \"Step 2:\", the paragraph(s) that followed it (the trailing\n// empty paragraph holding the _GoBack bookmark) are now directly after the\n// \"composer require...\" paragraph. Merge that following paragraph's content\n// up into the composer paragraph by deleting the paragraph mark between them,\n// i.e. extend the composer paragraph's range through the following\n// paragraph and collapse them into one.\nconst nextPara = composerPara.getNextOrNullObject();\nnextPara.load(\"isNullObject\");\nawait context.sync();\n\nif (!nextPara.isNullObject) {\n  // Range spanning from the end of the composer paragraph to the end of the\n  // following paragraph; deleting it removes the paragraph break that\n  // separates them, joining the two paragraphs into one (the following\n  // paragraph's trailing content, e.g. the bookmark, survives in the merged\n  // paragraph since only the paragraph mark itself is removed).\n  const composerEnd = composerPara.getRange(\"End\");\n  const nextEnd = nextPara.getRange(\"End\");\n  const joinRange = composerEnd.expandTo(nextEnd);\n  joinRange.delete();\n  await context.sync();\n}\n", "ps1": "# Remove the \"Step 2:\" paragraph and merge the trailing (empty, bookmark-only)\n# paragraph into the preceding \"composer require intervention/image\"\n# paragraph -- i.e. the \"Step 2:\" paragraph mark/text disappear entirely and\n# the _GoBack bookmark ends up living inside the \"composer require...\"\n# paragraph.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: delete the \"Step 2:\" paragraph (its text and paragraph mark). ---\n$stepTwoRange = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Trim() -eq \"Step 2:\") {\n        $stepTwoRange = $p.Range\n        break\n    }\n}\nif ($stepTwoRange -ne $null) {\n    $stepTwoRange.Delete()\n}\n\n# --- Step 2: merge the \"composer require intervention/image\" paragraph with\n# the paragraph that now immediately follows it (the bookmark-only empty\n# paragraph), by deleting the paragraph mark that separates them. ---\n$composerPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Trim() -eq \"composer require intervention/image\") {\n        $composerPara = $p\n        break\n    }\n}\nif ($composerPara -ne $null) {\n    $paraEnd = $composerPara.Range.End\n    $markRange = $d.Range($paraEnd - 1, $paraEnd)\n    $markRange.Delete()\n}\n"}
